{"js": "// Apply the \"Gods of Gold Infinireels\" copy-edit:\n//  - refresh the title (drop \"for\"/\"2021\") in the H1 and the bold recap line\n//  - rewrite three \"What we like\" bullets and one \"What we don't like\" bullet\n//  - rewrite the closing italic summary line\n// All target strings are used verbatim exactly where they should be replaced,\n// so a simple body-wide search/replace per string is sufficient (the title\n// string happens to repeat twice and both spots get the same replacement).\n\nconst replacements = [\n  [\n    \"Play Gods of Gold Infinireels for Free - Review 2021\",\n    \"Play Gods of Gold Infinireels Free - Review\",\n  ],\n  [\n    \"Innovative and unique infinite reel payment system\",\n    \"Innovative infinite reel payment system\",\n  ],\n  [\n    \"Free spins bonus and two reel modifiers with guaranteed wins\",\n    \"Exciting bonus features and modifiers\",\n  ],\n  [\n    \"Optimized for mobile devices\",\n    \"High volatility with big prizes\",\n  ],\n  [\n    \"High volatility may not be suitable for all players\",\n    \"Limited bonus features\",\n  ],\n  [\n    \"Take a spin and play Gods of Gold Infinireels for free. Read our latest review of this innovative NetEnt online slot with unique infinite reel payment system.\",\n    \"Discover the exciting gameplay and impressive graphics of Gods of Gold Infinireels in this review. Play for free now!\",\n  ],\n];\n\nfor (const [searchText, newText] of replacements) {\n  const results = context.document.body.search(searchText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the \"Gods of Gold Infinireels\" copy-edit:\n#  - refresh the title (drop \"for\"/\"2021\") in the H1 and the bold recap line\n#  - rewrite three \"What we like\" bullets and one \"What we don't like\" bullet\n#  - rewrite the closing italic summary line\n# The title string repeats twice (H1 + bold recap) and both get the same\n# replacement, so ReplaceAll over the whole document body works for every\n# pair below.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"Play Gods of Gold Infinireels for Free - Review 2021\", \"Play Gods of Gold Infinireels Free - Review\"),\n    @(\"Innovative and unique infinite reel payment system\", \"Innovative infinite reel payment system\"),\n    @(\"Free spins bonus and two reel modifiers with guaranteed wins\", \"Exciting bonus features and modifiers\"),\n    @(\"Optimized for mobile devices\", \"High volatility with big prizes\"),\n    @(\"High volatility may not be suitable for all players\", \"Limited bonus features\"),\n    @(\"Take a spin and play Gods of Gold Infinireels for free. Read our latest review of this innovative NetEnt online slot with unique infinite reel payment system.\", \"Discover the exciting gameplay and impressive graphics of Gods of Gold Infinireels in this review. Play for free now!\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute([ref]$findText, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]$replaceText, [ref]2) | Out-Null\n}\n"}
